# "Generate Report for Archive"
#
# 1. Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2. Narrow the "Status" column(s): Overview columns E & F, and column C
#    on both the zh-cn and de-de sheets (17.21598... -> 13.41018... chars,
#    i.e. a "character width" of 12.5, the closest this engine's
#    integer-pixel column-width model can represent).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- 1. Update the status cells ---
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"

$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"

# --- 2. Narrow the Status columns ---
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

$ws2.Columns.Item(3).ColumnWidth = 12.5

$ws3.Columns.Item(3).ColumnWidth = 12.5
